$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column D a bit (28.850625 -> 24.790625)
$ws.Columns.Item(4).ColumnWidth = 24

# Row 12 : carrier code + carrier name, and re-timestamp the survey date
$ws.Range("A12").Value = "TEST01"
$ws.Range("D12").Value = "100-LO CICERO ANTONINO"
$ws.Range("C12").Value = Get-Date -Year 2021 -Month 7 -Day 30 -Hour 11 -Minute 47 -Second 51

# Row 13 : carrier code + carrier name, and re-timestamp the survey date
$ws.Range("A13").Value = "TEST02"
$ws.Range("D13").Value = "102-LOGISTICA NIEDDU "
$ws.Range("C13").Value = Get-Date -Year 2021 -Month 7 -Day 30 -Hour 11 -Minute 48 -Second 3

# Row 14 : carrier code, carrier name now matches row 13's carrier, re-timestamp the survey date
$ws.Range("A14").Value = "TEST03"
$ws.Range("D14").Value = "102-LOGISTICA NIEDDU "
$ws.Range("C14").Value = Get-Date -Year 2021 -Month 7 -Day 30 -Hour 11 -Minute 49 -Second 1
